$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3.0
$ws.Range("G2").Value = 162.399297
$ws.Range("H2").Value = 487.197891
$ws.Range("I2").Value = 0.3910371682630009
$ws.Range("J2").Value = 0.3910371682630009
$ws.Range("K2").Value = 3.0
$ws.Range("M2").Value = 0.8019210000000001
$ws.Range("N2").Value = 2.405763
$ws.Range("O2").Value = 0.1587003260603691
$ws.Range("P2").Value = 0.1587003260603691
$ws.Range("Q2").Value = 130.231406649537
$ws.Range("R2").Value = 1172.082659845833
$ws.Range("S2").Value = 0.06205772610506165
$ws.Range("T2").Value = 0.06205772610506165
$ws.Range("E3").Value = 3.0
$ws.Range("G3").Value = 162.399297
$ws.Range("H3").Value = 487.197891
$ws.Range("I3").Value = 0.3910371682630009
$ws.Range("J3").Value = 0.3910371682630009
$ws.Range("K3").Value = 3.0
$ws.Range("M3").Value = 1.478072333333333
$ws.Range("N3").Value = 4.434217
$ws.Range("O3").Value = 0.2925108099685761
$ws.Range("P3").Value = 0.2925108099685761
$ws.Range("Q3").Value = 240.037907848483
$ws.Range("R3").Value = 2160.341170636347
$ws.Range("S3").Value = 0.1143825988164288
$ws.Range("T3").Value = 0.1143825988164288
$ws.Range("E4").Value = 3.0
$ws.Range("G4").Value = 162.399297
$ws.Range("H4").Value = 487.197891
$ws.Range("I4").Value = 0.3910371682630009
$ws.Range("J4").Value = 0.3910371682630009
$ws.Range("K4").Value = 3.0
$ws.Range("M4").Value = 0.2551363333333334
$ws.Range("N4").Value = 0.765409
$ws.Range("O4").Value = 0.05049153132272008
$ws.Range("P4").Value = 0.05049153132272008
$ws.Range("Q4").Value = 41.433961172491
$ws.Range("R4").Value = 372.905650552419
$ws.Range("S4").Value = 0.01974406542969907
$ws.Range("T4").Value = 0.01974406542969907
$ws.Range("E5").Value = 3.0
$ws.Range("G5").Value = 162.399297
$ws.Range("H5").Value = 487.197891
$ws.Range("I5").Value = 0.3910371682630009
$ws.Range("J5").Value = 0.3910371682630009
$ws.Range("K5").Value = 3.0
$ws.Range("M5").Value = 2.517922333333333
$ws.Range("N5").Value = 7.553767
$ws.Range("O5").Value = 0.4982973326483348
$ws.Range("P5").Value = 0.4982973326483348
$ws.Range("Q5").Value = 408.908816833933
$ws.Range("R5").Value = 3680.179351505396
$ws.Range("S5").Value = 0.1948527779118114
$ws.Range("T5").Value = 0.1948527779118114
$ws.Range("E6").Value = 3.0
$ws.Range("G6").Value = 65.41736466666667
$ws.Range("H6").Value = 196.252094
$ws.Range("I6").Value = 0.1575168212364948
$ws.Range("J6").Value = 0.1575168212364948
$ws.Range("K6").Value = 3.0
$ws.Range("M6").Value = 0.8019210000000001
$ws.Range("N6").Value = 2.405763
$ws.Range("O6").Value = 0.1587003260603691
$ws.Range("P6").Value = 0.1587003260603691
$ws.Range("Q6").Value = 52.45955849085801
$ws.Range("R6").Value = 472.1360264177221
$ws.Range("S6").Value = 0.0249979708902246
$ws.Range("T6").Value = 0.0249979708902246
$ws.Range("E7").Value = 3.0
$ws.Range("G7").Value = 65.41736466666667
$ws.Range("H7").Value = 196.252094
$ws.Range("I7").Value = 0.1575168212364948
$ws.Range("J7").Value = 0.1575168212364948
$ws.Range("K7").Value = 3.0
$ws.Range("M7").Value = 1.478072333333333
$ws.Range("N7").Value = 4.434217
$ws.Range("O7").Value = 0.2925108099685761
$ws.Range("P7").Value = 0.2925108099685761
$ws.Range("Q7").Value = 96.69159683337757
$ws.Range("R7").Value = 870.2243715003981
$ws.Range("S7").Value = 0.04607537296356251
$ws.Range("T7").Value = 0.04607537296356251
$ws.Range("E8").Value = 3.0
$ws.Range("G8").Value = 65.41736466666667
$ws.Range("H8").Value = 196.252094
$ws.Range("I8").Value = 0.1575168212364948
$ws.Range("J8").Value = 0.1575168212364948
$ws.Range("K8").Value = 3.0
$ws.Range("M8").Value = 0.2551363333333334
$ws.Range("N8").Value = 0.765409
$ws.Range("O8").Value = 0.05049153132272008
$ws.Range("P8").Value = 0.05049153132272008
$ws.Range("Q8").Value = 16.69034655738289
$ws.Range("R8").Value = 150.213119016446
$ws.Range("S8").Value = 0.007953265513317778
$ws.Range("T8").Value = 0.007953265513317778
$ws.Range("E9").Value = 3.0
$ws.Range("G9").Value = 65.41736466666667
$ws.Range("H9").Value = 196.252094
$ws.Range("I9").Value = 0.1575168212364948
$ws.Range("J9").Value = 0.1575168212364948
$ws.Range("K9").Value = 3.0
$ws.Range("M9").Value = 2.517922333333333
$ws.Range("N9").Value = 7.553767
$ws.Range("O9").Value = 0.4982973326483348
$ws.Range("P9").Value = 0.4982973326483348
$ws.Range("Q9").Value = 164.7158434820109
$ws.Range("R9").Value = 1482.442591338098
$ws.Range("S9").Value = 0.07849021186938994
$ws.Range("T9").Value = 0.07849021186938994
$ws.Range("E10").Value = 3.0
$ws.Range("G10").Value = 126.3069433333333
$ws.Range("H10").Value = 378.92083
$ws.Range("I10").Value = 0.3041313008456065
$ws.Range("J10").Value = 0.3041313008456065
$ws.Range("K10").Value = 3.0
$ws.Range("M10").Value = 0.8019210000000001
$ws.Range("N10").Value = 2.405763
$ws.Range("O10").Value = 0.1587003260603691
$ws.Range("P10").Value = 0.1587003260603691
$ws.Range("Q10").Value = 101.28819030481
$ws.Range("R10").Value = 911.5937127432902
$ws.Range("S10").Value = 0.04826573660936195
$ws.Range("T10").Value = 0.04826573660936195
$ws.Range("E11").Value = 3.0
$ws.Range("G11").Value = 126.3069433333333
$ws.Range("H11").Value = 378.92083
$ws.Range("I11").Value = 0.3041313008456065
$ws.Range("J11").Value = 0.3041313008456065
$ws.Range("K11").Value = 3.0
$ws.Range("M11").Value = 1.478072333333333
$ws.Range("N11").Value = 4.434217
$ws.Range("O11").Value = 0.2925108099685761
$ws.Range("P11").Value = 0.2925108099685761
$ws.Range("Q11").Value = 186.6907984489011
$ws.Range("R11").Value = 1680.21718604011
$ws.Range("S11").Value = 0.08896169314714504
$ws.Range("T11").Value = 0.08896169314714505
$ws.Range("E12").Value = 3.0
$ws.Range("G12").Value = 126.3069433333333
$ws.Range("H12").Value = 378.92083
$ws.Range("I12").Value = 0.3041313008456065
$ws.Range("J12").Value = 0.3041313008456065
$ws.Range("K12").Value = 3.0
$ws.Range("M12").Value = 0.2551363333333334
$ws.Range("N12").Value = 0.765409
$ws.Range("O12").Value = 0.05049153132272008
$ws.Range("P12").Value = 0.05049153132272008
$ws.Range("Q12").Value = 32.22549039660778
$ws.Range("R12").Value = 290.02941356947
$ws.Range("S12").Value = 0.01535605510286554
$ws.Range("T12").Value = 0.01535605510286554
$ws.Range("E13").Value = 3.0
$ws.Range("G13").Value = 126.3069433333333
$ws.Range("H13").Value = 378.92083
$ws.Range("I13").Value = 0.3041313008456065
$ws.Range("J13").Value = 0.3041313008456065
$ws.Range("K13").Value = 3.0
$ws.Range("M13").Value = 2.517922333333333
$ws.Range("N13").Value = 7.553767
$ws.Range("O13").Value = 0.4982973326483348
$ws.Range("P13").Value = 0.4982973326483348
$ws.Range("Q13").Value = 318.0310734740678
$ws.Range("R13").Value = 2862.27966126661
$ws.Range("S13").Value = 0.1515478159862339
$ws.Range("T13").Value = 0.151547815986234
$ws.Range("E14").Value = 3.0
$ws.Range("G14").Value = 61.180387
$ws.Range("H14").Value = 183.541161
$ws.Range("I14").Value = 0.1473147096548978
$ws.Range("J14").Value = 0.1473147096548978
$ws.Range("K14").Value = 3.0
$ws.Range("M14").Value = 0.8019210000000001
$ws.Range("N14").Value = 2.405763
$ws.Range("O14").Value = 0.1587003260603691
$ws.Range("P14").Value = 0.1587003260603691
$ws.Range("Q14").Value = 49.061837123427
$ws.Range("R14").Value = 441.556534110843
$ws.Range("S14").Value = 0.02337889245572088
$ws.Range("T14").Value = 0.02337889245572088
$ws.Range("E15").Value = 3.0
$ws.Range("G15").Value = 61.180387
$ws.Range("H15").Value = 183.541161
$ws.Range("I15").Value = 0.1473147096548978
$ws.Range("J15").Value = 0.1473147096548978
$ws.Range("K15").Value = 3.0
$ws.Range("M15").Value = 1.478072333333333
$ws.Range("N15").Value = 4.434217
$ws.Range("O15").Value = 0.2925108099685761
$ws.Range("P15").Value = 0.2925108099685761
$ws.Range("Q15").Value = 90.42903736732633
$ws.Range("R15").Value = 813.861336305937
$ws.Range("S15").Value = 0.04309114504143977
$ws.Range("T15").Value = 0.04309114504143978
$ws.Range("E16").Value = 3.0
$ws.Range("G16").Value = 61.180387
$ws.Range("H16").Value = 183.541161
$ws.Range("I16").Value = 0.1473147096548978
$ws.Range("J16").Value = 0.1473147096548978
$ws.Range("K16").Value = 3.0
$ws.Range("M16").Value = 0.2551363333333334
$ws.Range("N16").Value = 0.765409
$ws.Range("O16").Value = 0.05049153132272008
$ws.Range("P16").Value = 0.05049153132272008
$ws.Range("Q16").Value = 15.60933961109433
$ws.Range("R16").Value = 140.484056499849
$ws.Range("S16").Value = 0.007438145276837686
$ws.Range("T16").Value = 0.007438145276837686
$ws.Range("E17").Value = 3.0
$ws.Range("G17").Value = 61.180387
$ws.Range("H17").Value = 183.541161
$ws.Range("I17").Value = 0.1473147096548978
$ws.Range("J17").Value = 0.1473147096548978
$ws.Range("K17").Value = 3.0
$ws.Range("M17").Value = 2.517922333333333
$ws.Range("N17").Value = 7.553767
$ws.Range("O17").Value = 0.4982973326483348
$ws.Range("P17").Value = 0.4982973326483348
$ws.Range("Q17").Value = 154.0474627892763
$ws.Range("R17").Value = 1386.427165103487
$ws.Range("S17").Value = 0.07340652688089945
$ws.Range("T17").Value = 0.07340652688089946
